$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 962.08
$ws.Range("I33").Value = 692.2105
$ws.Range("J33").Value = 1816.6666
$ws.Range("K33").Value = 692.2105
$ws.Range("L33").Value = 1816.6666
$ws.Range("M33").Value = -463.2105
$ws.Range("N33").Value = -2274.6666
$ws.Range("H64").Value = 3033.3333
$ws.Range("I64").Value = 3033.3333
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3033.3333
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -2785.3333
$ws.Range("H67").Value = 3033.3333
$ws.Range("I67").Value = 3033.3333
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3033.3333
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -2175.3333
$ws.Range("H86").Value = 140015200
$ws.Range("I86").Value = 175018160
$ws.Range("J86").Value = 3336
$ws.Range("K86").Value = 175018160
$ws.Range("L86").Value = 3336
$ws.Range("M86").Value = -175017037
$ws.Range("N86").Value = -5582
$ws.Range("H89").Value = 140015200
$ws.Range("I89").Value = 175018160
$ws.Range("J89").Value = 3336
$ws.Range("K89").Value = 875090800
$ws.Range("L89").Value = 16680
$ws.Range("M89").Value = -875085184
$ws.Range("N89").Value = -27912
$ws.Range("H105").Value = 59000
$ws.Range("J105").Value = 59000
$ws.Range("L105").Value = 59000
$ws.Range("N105").Value = -65988
$ws.Range("H138").Value = 2243.3872
$ws.Range("I138").Value = 2679.4666
$ws.Range("J138").Value = 2104.2126
$ws.Range("K138").Value = 8038.399800000001
$ws.Range("L138").Value = 6312.6378
$ws.Range("M138").Value = -2898.399800000001
$ws.Range("N138").Value = -16592.6378

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 6779
$ws.Range("I28").Value = 3461.7144
$ws.Range("J28").Value = 30000
$ws.Range("K28").Value = 3461.7144
$ws.Range("L28").Value = 30000
$ws.Range("M28").Value = -3269.7144
$ws.Range("N28").Value = -30384
$ws.Range("H32").Value = 424983.03
$ws.Range("I32").Value = 507882.66
$ws.Range("K32").Value = 507882.66
$ws.Range("M32").Value = -507595.66
$ws.Range("H61").Value = 3058.303
$ws.Range("I61").Value = 2489.0588
$ws.Range("J61").Value = 3663.125
$ws.Range("K61").Value = 2489.0588
$ws.Range("L61").Value = 3663.125
$ws.Range("M61").Value = -2277.0588
$ws.Range("N61").Value = -4087.125
$ws.Range("H88").Value = 2754.1667
$ws.Range("I88").Value = 2703.6
$ws.Range("K88").Value = 2703.6
$ws.Range("M88").Value = -2297.6
$ws.Range("H91").Value = 2754.1667
$ws.Range("I91").Value = 2703.6
$ws.Range("K91").Value = 2703.6
$ws.Range("M91").Value = -1299.6
$ws.Range("H92").Value = 59704
$ws.Range("J92").Value = 59704
$ws.Range("L92").Value = 59704
$ws.Range("N92").Value = -64696
$ws.Range("H96").Value = 43922
$ws.Range("J96").Value = 43922
$ws.Range("L96").Value = 43922
$ws.Range("N96").Value = -49414
$ws.Range("H99").Value = 6779
$ws.Range("I99").Value = 3461.7144
$ws.Range("J99").Value = 30000
$ws.Range("K99").Value = 3461.7144
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -466.7143999999998
$ws.Range("N99").Value = -35990
$ws.Range("H104").Value = 73333.336
$ws.Range("J104").Value = 73333.336
$ws.Range("L104").Value = 73333.336
$ws.Range("N104").Value = -80321.336
$ws.Range("H136").Value = 3058.303
$ws.Range("I136").Value = 2489.0588
$ws.Range("J136").Value = 3663.125
$ws.Range("K136").Value = 7467.176399999999
$ws.Range("L136").Value = 10989.375
$ws.Range("M136").Value = -4917.176399999999
$ws.Range("N136").Value = -16089.375

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 41668492
$ws.Range("I86").Value = 47620852
$ws.Range("J86").Value = 1969
$ws.Range("K86").Value = 47620852
$ws.Range("L86").Value = 1969
$ws.Range("M86").Value = -47619729
$ws.Range("N86").Value = -4215
$ws.Range("H89").Value = 41668492
$ws.Range("I89").Value = 47620852
$ws.Range("J89").Value = 1969
$ws.Range("K89").Value = 238104260
$ws.Range("L89").Value = 9845
$ws.Range("M89").Value = -238098644
$ws.Range("N89").Value = -21077
$ws.Range("H97").Value = 27499
$ws.Range("J97").Value = 34998
$ws.Range("L97").Value = 34998
$ws.Range("N97").Value = -36980
$ws.Range("H107").Value = 168500
$ws.Range("I107").Value = 168500
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 168500
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -166580

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5908.58
$ws.Range("I31").Value = 1521.7778
$ws.Range("J31").Value = 11058.305
$ws.Range("K31").Value = 1521.7778
$ws.Range("L31").Value = 11058.305
$ws.Range("M31").Value = -1226.7778
$ws.Range("N31").Value = -11648.305
$ws.Range("H34").Value = 5908.58
$ws.Range("I34").Value = 1521.7778
$ws.Range("J34").Value = 11058.305
$ws.Range("K34").Value = 1521.7778
$ws.Range("L34").Value = 11058.305
$ws.Range("M34").Value = -1319.7778
$ws.Range("N34").Value = -11462.305
$ws.Range("H58").Value = 1136.475
$ws.Range("I58").Value = 853.4815
$ws.Range("J58").Value = 1724.2307
$ws.Range("K58").Value = 853.4815
$ws.Range("L58").Value = 1724.2307
$ws.Range("M58").Value = -650.4815
$ws.Range("N58").Value = -2130.2307
$ws.Range("H136").Value = 1136.475
$ws.Range("I136").Value = 853.4815
$ws.Range("J136").Value = 1724.2307
$ws.Range("K136").Value = 2560.4445
$ws.Range("L136").Value = 5172.6921
$ws.Range("M136").Value = -10.44450000000006
$ws.Range("N136").Value = -10272.6921

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 769.34485
$ws.Range("I131").Value = 230.90909
$ws.Range("J131").Value = 1098.3889
$ws.Range("K131").Value = 692.72727
$ws.Range("L131").Value = 3295.1667
$ws.Range("M131").Value = 4347.27273
$ws.Range("N131").Value = -13375.1667
$ws.Range("H140").Value = 1243.25
$ws.Range("I140").Value = 1050.75
$ws.Range("K140").Value = 3152.25
$ws.Range("M140").Value = 2027.75

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 20000
$ws.Range("J96").Value = 20000
$ws.Range("L96").Value = 20000
$ws.Range("N96").Value = -25492
$ws.Range("H102").Value = 1432.6
$ws.Range("I102").Value = 1252.7778
$ws.Range("J102").Value = 1702.3334
$ws.Range("K102").Value = 1252.7778
$ws.Range("L102").Value = 1702.3334
$ws.Range("M102").Value = 369.2221999999999
$ws.Range("N102").Value = -4946.3334

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4387692
$ws.Range("I136").Value = 1299.55
$ws.Range("J136").Value = 9261462
$ws.Range("K136").Value = 3898.65
$ws.Range("L136").Value = 27784386
$ws.Range("M136").Value = -1348.65
$ws.Range("N136").Value = -27789486

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3625492.5
$ws.Range("I132").Value = 3344.0715
$ws.Range("J132").Value = 5210182.5
$ws.Range("K132").Value = 10032.2145
$ws.Range("L132").Value = 15630547.5
$ws.Range("M132").Value = -7502.2145
$ws.Range("N132").Value = -15635607.5

Write-Output "Applied all Anima_Profits updates"